$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.197.74"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "3.164.28"
$ws.Range("E3").Value = "  +3.23%  "

$ws.Range("E4").Value = "  +0.09%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.12"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +3.50%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.52"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +5.43%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.161.22"
$ws.Range("E8").Value = "  +3.07%  "

$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("E10").Value = "  +5.43%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("E12").Value = "  +3.95%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +15.59%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.36"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +5.86%  "

$ws.Range("D15").Value = "3.685.31"
$ws.Range("E15").Value = "  +3.35%  "

$ws.Range("D16").Value = "65.268.23"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17").Value = "3.164.68"
$ws.Range("E17").Value = "  +3.57%  "

$ws.Range("E18").Value = "  +5.35%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.76"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +4.69%  "

$ws.Range("E21").Value = "  +3.91%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +4.83%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.34"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +5.42%  "

$ws.Range("E24").Value = "  +3.38%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.86"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +2.53%  "

$ws.Range("E26").Value = "  -0.01%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.01"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +10.61%  "

$ws.Range("E28").Value = "  +4.48%  "

$ws.Range("E29").Value = "  +7.15%  "

$ws.Range("E30").Value = "  +14.06%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.83"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +4.54%  "

$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("E33").Value = "  +3.93%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.32"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +10.65%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.57"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +5.50%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.29"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -0.07%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0908"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +11.20%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "473.17"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +6.48%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0424"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +2.67%  "

$ws.Range("E40").Value = "  +9.39%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.68"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +3.89%  "

$ws.Range("D42").Value = "3.066.78"
$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("E43").Value = "  +1.51%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.46"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +10.49%  "

$ws.Range("E45").Value = "  +4.09%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.93"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +3.83%  "

$ws.Range("D47").Value = "0.0₃0597"
$ws.Range("E47").Value = "  +15.57%  "

$ws.Range("E49").Value = "  +0.77%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +6.38%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.88"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +1.91%  "
